$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data to row 8 (Hours / Completed)
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "Load data from Parse.com and display them on the ArticlesTableViewController"

# Widen column B to fit the new content
$ws.Columns.Item(2).ColumnWidth = 20

# Move the active selection to C9
$ws.Range("C9").Select()
